# ReimbUnitsCreationScenarios.xlsx - "Reimb Form Object classes" edit
#
# 1. The second test-data row (AutoUnit2 / "Reimb Unit Form auto2") is removed
#    entirely - delete worksheet row 3, which shifts every row below it up by one.
# 2. The remaining row's test-name cell (D2) is renamed from
#    "Create Custom flow1" to "Create Reimbursement flow1".
# 3. C2:E2 lose the custom style they had (10pt Arial) and fall back to the
#    default cell style already used by A2/B2 - copy A2's format onto C2:E2.
# 4. The saved selection moves from C8 to B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the remaining flow description before the row shift so the edit
# lands on the surviving data row.
$ws.Range("D2").Value = "Create Reimbursement flow1"

# Drop the second (AutoUnit2) test case row; rows below shift up.
$ws.Rows(3).Delete()

# C2:E2 pick up the plain style already used by A2/B2 instead of the bold
# style that only the removed row's siblings needed.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("C2:E2").PasteSpecial(-4122) | Out-Null

# Update the saved cursor position.
$ws.Range("B5").Select()
